$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.734869837760925
$ws.Range("B1").Value = 5.120697975158691
$ws.Range("C1").Value = 4.062240123748779
$ws.Range("D1").Value = 0.980292797088623
$ws.Range("E1").Value = 0.5987050533294678
